$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'mens compression basketball pants'
$ws.Range("A2").Value = 'mens basketball tights with knee pads'
$ws.Range("A3").Value = 'capri pants for boys'
$ws.Range("A4").Value = 'knee pads with leggings'
$ws.Range("A5").Value = 'compression pants for men'
$ws.Range("A6").Value = 'compression mens leggings'
$ws.Range("A7").Value = 'basketball protective gear'
$ws.Range("A8").Value = 'basketball spandex pants'
$ws.Range("A9").Value = 'compression tights girls'
$ws.Range("A10").Value = 'wrestling pants'
$ws.Range("A11").Value = 'mens compression capris'
$ws.Range("A12").Value = 'baseball sliding tights'
$ws.Range("A13").Value = 'boys basketball compression pants with knee pads'
$ws.Range("A14").Value = 'compression men capri'
$ws.Range("A15").Value = 'volleyball knee pads men'
$ws.Range("A16").Value = 'mens basketball'
$ws.Range("A17").Value = 'athletic leggings boys'
$ws.Range("A18").Value = 'knee pad weightlifting'
$ws.Range("A19").Value = 'basketball guide'
$ws.Range("A20").Value = 'basketball knee pads boys'
$ws.Range("A21").Value = 'youth padded compression pants basketball'
$ws.Range("A22").Value = 'leggings for men'
$ws.Range("A23").Value = 'cycling capris'
$ws.Range("A24").Value = 'padded knee pads for basketball'
$ws.Range("A25").Value = 'pant with knee pad'
$ws.Range("A26").Value = 'tight pants for men'
$ws.Range("A27").Value = 'basketball kneepads'
$ws.Range("A28").Value = 'thigh compression leggings'
$ws.Range("A29").Value = 'compression pants men pack'
$ws.Range("A30").Value = 'hiking knee pads'
$ws.Range("A31").Value = 'sliding pants'
$ws.Range("A32").Value = 'boys compression pants'
$ws.Range("A33").Value = 'cycling knee pads for men'
$ws.Range("A34").Value = 'youth boys knee pads for basketball'
$ws.Range("A35").Value = 'mens spandex leggings'
$ws.Range("A36").Value = 'compression knee padded'
$ws.Range("A37").Value = 'youth basketball kneepads'
$ws.Range("A38").Value = 'baseball equipment for boys'
$ws.Range("A39").Value = 'youth football tights boys'
$ws.Range("A40").Value = 'men compression workout pants'
$ws.Range("A41").Value = 'padded knee basketball'
$ws.Range("A42").Value = 'mens compression pants running'
$ws.Range("A43").Value = 'basketball youth leggings'
$ws.Range("A44").Value = 'compression basketball leggings youth'
$ws.Range("A45").Value = 'men compression running pants'
$ws.Range("A46").Value = 'wrestling knee pads youth'
$ws.Range("A47").Value = 'mens leggings sports'
$ws.Range("A48").Value = 'skin leggings men'
$ws.Range("A49").Value = 'youth football pants'
$ws.Range("A50").Value = 'mens leggings for sports'
$ws.Range("A51").Value = 'black boys compression pants'
$ws.Range("A52").Value = 'basketball team clothes'
$ws.Range("A53").Value = 'men compression tights'
$ws.Range("A54").Value = 'youth spandex leggings boys'
$ws.Range("A55").Value = 'baseball gear for men'
$ws.Range("A56").Value = 'sports knee pads'
$ws.Range("A57").Value = 'leggings men'
$ws.Range("A58").Value = 'pain in hip down leg'
$ws.Range("A59").Value = 'compression men tights'
$ws.Range("A60").Value = 'softball pants for men'
$ws.Range("A61").Value = 'youth basketball knee pads for boys'
$ws.Range("A62").Value = 'mens running tights'
$ws.Range("A63").Value = 'youth baseball gear'
$ws.Range("A64").Value = 'football pants adult xl'
$ws.Range("A65").Value = 'padded knee compression'
$ws.Range("A66").Value = 'compression pants men 3/4'
$ws.Range("A67").Value = 'sliding pants softball youth'
$ws.Range("A68").Value = 'hockey hip pads'
$ws.Range("A69").Value = 'compression knees'
$ws.Range("A70").Value = 'compression men'
$ws.Range("A71").Value = 'compression mens tights'
$ws.Range("A72").Value = 'compression tights men'
$ws.Range("A73").Value = 'legging for men'
$ws.Range("A74").Value = 'football knee pads adult'
$ws.Range("A75").Value = 'football pants with pads mens'
$ws.Range("A76").Value = 'mens basketball tights and leggings'
$ws.Range("A77").Value = 'knee protector'
$ws.Range("A78").Value = 'basketball apparel mens'
$ws.Range("A79").Value = 'adult knee pads'
$ws.Range("A80").Value = 'black softball pants girls'
$ws.Range("A81").Value = 'youth football pads for pants'
$ws.Range("A82").Value = 'youth black compression pants'
$ws.Range("A83").Value = '3/4 compression pants men'
$ws.Range("A84").Value = 'girls hiking pants'
$ws.Range("A85").Value = 'youth tights boys basketball'
$ws.Range("A86").Value = 'down pants men'
$ws.Range("A87").Value = 'basketball volleyball knee pads'
$ws.Range("A88").Value = 'legging for men sport'
$ws.Range("A89").Value = 'polyester hex mesh'
$ws.Range("A90").Value = 'basketball knee pad tights'
$ws.Range("A91").Value = 'capris tights'
$ws.Range("A92").Value = 'thigh compression pants'
$ws.Range("A93").Value = 'bjj pants men'
$ws.Range("A94").Value = 'wrestling pads'
$ws.Range("A95").Value = 'basketball athletic tights'
$ws.Range("A96").Value = 'wrestling knee pads youth 2 pack'
$ws.Range("A97").Value = 'compression for knees'
$ws.Range("A98").Value = 'youth football pants with pads black'
$ws.Range("A99").Value = 'girls knee pads'
$ws.Range("A100").Value = 'compression pants football'
